$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.812.09'
$ws.Range("D3").Value = '1.625.63'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.83'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.518'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.19'
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").Value = '1.856.14'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '1.658.88'
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.555'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.86'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = '27.825.20'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.79'
$ws.Range("E18").Value = '  -2.02%  '
$ws.Range("D19").Value = '0.0₃0717'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.56'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("E23").Value = '  -4.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.90'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.40'
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.17'
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '1.402.31'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.845'
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.995'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.57'
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("E44").Value = '  -1.45%  '
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("D46").Value = '1.765.60'
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.18'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("E51").Value = '  +0.05%  '
